# Colombia shapefile from R
# Updates the DSSAT input-file-extension callouts on slide 2 to reflect
# the real Colombia shapefile extensions produced by R ([.WTH], [.SOL],
# [.X]) and resizes their rectangles to fit the new label widths. Also
# joins the wrapped "Decision Support System for Agrotechnology Transfer"
# subtitle back into a single run.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# Shape.Width/.Height are points backed by a 32-bit float, so a plain
# EMU/12700 division can land one EMU short after the round-trip back to
# EMU on save; nudge by a hair so it rounds to the exact target EMU.
function Set-WidthEmu($shape, $emu) {
    $shape.Width = ($emu / 12700.0) + 0.00002
}

# ---------------------------------------------------------------
# [.wth files]  ->  [.WTH files]   (Rectangle 17)
# ---------------------------------------------------------------
$shpWth = $s.Shapes.Item(11)
Set-WidthEmu $shpWth 1316194
$trWth = $shpWth.TextFrame.TextRange
$subWth = $trWth.Characters(1, 6)
$subWth.Text = "[.WTH "

# ---------------------------------------------------------------
# [.soil files]  ->  [.SOL files]   (Rectangle 20)
# ---------------------------------------------------------------
$shpSol = $s.Shapes.Item(13)
Set-WidthEmu $shpSol 1226618
$trSol = $shpSol.TextFrame.TextRange
$subSol = $trSol.Characters(1, 7)
$subSol.Text = "[.SOL "

# ---------------------------------------------------------------
# [.x files]  ->  [.X files]   (Rectangle 23)
# ---------------------------------------------------------------
$shpX = $s.Shapes.Item(16)
Set-WidthEmu $shpX 990977
$trX = $shpX.TextFrame.TextRange
$subX = $trX.Characters(1, 4)
$subX.Text = "[.X "

# ---------------------------------------------------------------
# Merge "Decision Support System for " + "Agrotechnology" + " Transfer"
# back into one run (TextBox 8).
# ---------------------------------------------------------------
$shpTitle = $s.Shapes.Item(21)
$titleHeight = $shpTitle.Height
$trTitle = $shpTitle.TextFrame.TextRange
$prefix = $trTitle.Characters(1, 7)
$prefixText = $prefix.Text
$prefixSize = $prefix.Font.Size
$trTitle.Delete()
$trTitle.Text = $prefixText
$trTitle.Font.Size = $prefixSize
$null = $trTitle.InsertAfter("Decision Support System for Agrotechnology Transfer")
$suffix = $trTitle.Characters(8, 51)
$suffix.Font.Size = 24
# Re-merging the runs perturbs the spAutoFit height calculation; restore
# the original autofit height since the wrapped text content/size is
# unchanged (only the run split changed).
$shpTitle.Height = $titleHeight
